# Applies the "shift one new weekly record in" edit described by the diff:
# - Rows 96..196 take on the values (D, I, J, K, L, M, N, O, P, Q) that
#   previously belonged to rows 95..195 (i.e. every existing record is
#   pushed down by one row).
# - Row 95 is populated with a brand-new record (the most recent weekly
#   observation).
# - The sheet grows from A1:R195 to A1:R196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 95
$lastRow  = 195
$newLastRow = $lastRow + 1   # 196

# Columns that vary per-record and therefore need to shift down.
$varCols = 4,9,10,11,12,13,14,15,16,17   # D, I, J, K, L, M, N, O, P, Q

# Columns that are constant for every record in this sheet - needed so the
# brand new row (196) gets them too.
$constCols = 1,2,3,5,6,7,8,18            # A, B, C, E, F, G, H, R

# Preserve the date number format used by column D onto the new bottom row
# before writing a value into it, so Excel does not invent a brand new
# (differently formatted) style for it.
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

# 1) Snapshot the current values for the columns that move, for every
#    row in the block (95..195), before we start overwriting anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $varCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshot back out, shifted down by one row: what used to
#    be in row r now belongs in row r+1.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $destRow = $r + 1
    $rowVals = $snapshot[$r]
    foreach ($c in $varCols) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c]
    }
}

# Copy the constant columns down onto the brand-new last row too, since it
# did not exist before.
foreach ($c in $constCols) {
    $ws.Cells.Item($newLastRow, $c).Value = $ws.Cells.Item($lastRow, $c).Value()
}

# 3) Populate row 95 with the new weekly record.
$ws.Cells.Item(95, 4).Value  = 44484          # D - Fecha
$ws.Cells.Item(95, 9).Value  = "Primera"      # I - Calidad
$ws.Cells.Item(95, 10).Value = 1200           # J - Volumen
$ws.Cells.Item(95, 11).Value = 1000           # K - Precio minimo
$ws.Cells.Item(95, 12).Value = 1200           # L - Precio maximo
$ws.Cells.Item(95, 13).Value = 1100           # M - Precio promedio ponderado
$ws.Cells.Item(95, 14).Value = "$/paquete 5 unidades"   # N - Unidad de comercializacion
$ws.Cells.Item(95, 15).Value = "Región del Maule"       # O - Origen
$ws.Cells.Item(95, 16).Value = 220            # P - Precio $/Kg
$ws.Cells.Item(95, 17).Value = 5              # Q - Kg o Unidades
